$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.892.69"
$ws.Range("E2").Value = "  -0.97%  "
$ws.Range("D3").Value = "1.868.05"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.96"
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5092"
$ws.Range("E7").Value = "  -1.48%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3657"
$ws.Range("E8").Value = "  -2.69%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07175"
$ws.Range("E9").Value = "  +0.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8897"
$ws.Range("E10").Value = "  -0.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.58"
$ws.Range("E11").Value = "  -0.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07484"
$ws.Range("E12").Value = "  -0.67%  "
$ws.Range("D13").Value = "1.871.63"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.45"
$ws.Range("E14").Value = "  +5.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.219"
$ws.Range("E15").Value = "  -1.61%  "
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008489"
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("E18").Value = "  +0.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("D20").Value = "26.942.71"
$ws.Range("E20").Value = "  -0.90%  "
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("D22").Value = "2.108.87"
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("E23").Value = "  -1.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.73"
$ws.Range("E25").Value = "  +1.27%  "
$ws.Range("E26").Value = "  -3.21%  "
$ws.Range("E27").Value = "  -0.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.084"
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.40"
$ws.Range("E29").Value = "  +0.44%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.679"
$ws.Range("E30").Value = "  +0.38%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.701"
$ws.Range("E31").Value = "  +0.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09100"
$ws.Range("E32").Value = "  -1.58%  "
$ws.Range("E33").Value = "  -1.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7475"
$ws.Range("E34").Value = "  +3.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.973"
$ws.Range("E35").Value = "  -3.53%  "
$ws.Range("E36").Value = "  -0.94%  "
$ws.Range("E37").Value = "  +3.78%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.504"
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01986"
$ws.Range("E39").Value = "  -2.24%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5564"
$ws.Range("E40").Value = "  +5.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.071"
$ws.Range("E41").Value = "  -0.64%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.579"
$ws.Range("E42").Value = "  +1.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "115.44"
$ws.Range("E43").Value = "  -1.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.576"
$ws.Range("E44").Value = "  +3.21%  "
$ws.Range("E45").Value = "  +0.89%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4759"
$ws.Range("E46").Value = "  +3.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.000"
$ws.Range("E47").Value = "  +0.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.12"
$ws.Range("E48").Value = "  +1.51%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "37.04"
$ws.Range("E49").Value = "  +0.97%  "
$ws.Range("E50").Value = "  -0.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "62.92"
$ws.Range("E51").Value = "  -1.09%  "
